# Update the worksheet date header and the 25 division problems/answers
# in the practice table to the next day's values.
$d = $word.ActiveDocument

# Header date
$d.Content.Find.Execute("2024-11-25 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-26 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("405÷5=81, 0", $true, $false, $false, $false, $false, $true, 1, $false, "711÷6=118, 3", 2) | Out-Null
$d.Content.Find.Execute("730÷4=182, 2", $true, $false, $false, $false, $false, $true, 1, $false, "174÷7=24, 6", 2) | Out-Null
$d.Content.Find.Execute("863÷3=287, 2", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=251, 2", 2) | Out-Null
$d.Content.Find.Execute("632÷4=158, 0", $true, $false, $false, $false, $false, $true, 1, $false, "647÷6=107, 5", 2) | Out-Null
$d.Content.Find.Execute("633÷2=316, 1", $true, $false, $false, $false, $false, $true, 1, $false, "333÷8=41, 5", 2) | Out-Null
$d.Content.Find.Execute("239÷4=59, 3", $true, $false, $false, $false, $false, $true, 1, $false, "599÷5=119, 4", 2) | Out-Null
$d.Content.Find.Execute("258÷2=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "847÷3=282, 1", 2) | Out-Null
$d.Content.Find.Execute("627÷7=89, 4", $true, $false, $false, $false, $false, $true, 1, $false, "978÷2=489, 0", 2) | Out-Null
$d.Content.Find.Execute("841÷4=210, 1", $true, $false, $false, $false, $false, $true, 1, $false, "840÷5=168, 0", 2) | Out-Null
$d.Content.Find.Execute("927÷3=309, 0", $true, $false, $false, $false, $false, $true, 1, $false, "563÷8=70, 3", 2) | Out-Null
$d.Content.Find.Execute("165÷9=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "977÷3=325, 2", 2) | Out-Null
$d.Content.Find.Execute("146÷6=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "766÷8=95, 6", 2) | Out-Null
$d.Content.Find.Execute("488÷5=97, 3", $true, $false, $false, $false, $false, $true, 1, $false, "351÷5=70, 1", 2) | Out-Null
$d.Content.Find.Execute("997÷5=199, 2", $true, $false, $false, $false, $false, $true, 1, $false, "292÷3=97, 1", 2) | Out-Null
$d.Content.Find.Execute("794÷2=397, 0", $true, $false, $false, $false, $false, $true, 1, $false, "754÷6=125, 4", 2) | Out-Null
$d.Content.Find.Execute("744÷6=124, 0", $true, $false, $false, $false, $false, $true, 1, $false, "758÷4=189, 2", 2) | Out-Null
$d.Content.Find.Execute("123÷6=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "976÷9=108, 4", 2) | Out-Null
$d.Content.Find.Execute("997÷8=124, 5", $true, $false, $false, $false, $false, $true, 1, $false, "629÷7=89, 6", 2) | Out-Null
$d.Content.Find.Execute("912÷6=152, 0", $true, $false, $false, $false, $false, $true, 1, $false, "410÷9=45, 5", 2) | Out-Null
$d.Content.Find.Execute("146÷9=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "530÷7=75, 5", 2) | Out-Null
$d.Content.Find.Execute("651÷8=81, 3", $true, $false, $false, $false, $false, $true, 1, $false, "496÷4=124, 0", 2) | Out-Null
$d.Content.Find.Execute("883÷7=126, 1", $true, $false, $false, $false, $false, $true, 1, $false, "442÷5=88, 2", 2) | Out-Null
$d.Content.Find.Execute("957÷4=239, 1", $true, $false, $false, $false, $false, $true, 1, $false, "779÷8=97, 3", 2) | Out-Null
$d.Content.Find.Execute("618÷9=68, 6", $true, $false, $false, $false, $false, $true, 1, $false, "840÷9=93, 3", 2) | Out-Null
$d.Content.Find.Execute("582÷2=291, 0", $true, $false, $false, $false, $false, $true, 1, $false, "489÷8=61, 1", 2) | Out-Null
